$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert the two new rows (shifts old rows 5-10 down to 6-11, then 8-11 down to 9-12)
$ws.Rows("5").Insert()
$ws.Rows("8").Insert()

# 2) Update B2 (Prod UserName) to match C2's email
$ws.Range("B2").Value = "avayugundla@helenoftroy.com"

# 3) Populate the new rows. Order chosen to reproduce the shared-string allocation
#    order in the target file (Simple_product_10qty, then "12", then
#    Configurable Product_10qty).
$ws.Range("A8").Value = "Simple_product_10qty"
$ws.Range("V5").Value = "'12"
$ws.Range("A5").Value = "Configurable Product_10qty"

$ws.Range("U8").Value = "Hair Makeup - Turquoise"
$ws.Range("V8").Value = "'12"
$ws.Range("U5").Value = "Curl Defining Styling Soufflé (Mini 2 fl oz.)"

# 4) Fix the active selection / scroll position
$ws.Range("C3").Select()

# 5) Rebuild hyperlinks for the final row positions (the engine only supports
#    deleting/re-adding the whole collection, not a single link at a time).
$ws.Cells.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Lotuswave@123") | Out-Null
$ws.Range("E2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:vnarra@helenoftroy.com") | Out-Null
$ws.Range("B2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:avayugundla@helenoftroy.com") | Out-Null
$ws.Range("C2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Lotuswave@123") | Out-Null
$ws.Range("D2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("K9"), "mailto:vnarra@helenoftroy.com") | Out-Null
$ws.Range("K9").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:avayugundla@helenoftroy.com") | Out-Null
$ws.Range("B11").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:Lotuswave@123") | Out-Null
$ws.Range("D11").Style = "Hyperlink"
